$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Note: Price column (D) values are stored as plain text (e.g. thousands are
# separated with '.' like "41.919.17", and some values have significant
# trailing zeros like "99.00"). A leading apostrophe is used where the new
# value would otherwise be auto-parsed by Excel as a number, so the text is
# preserved exactly as authored.
$ws.Range('D2').Value = '41.919.17'
$ws.Range('E2').Value = '  +5.75%  '
$ws.Range('D3').Value = '2.233.43'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''231.71'
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').Value = '''61.69'
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.402'
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').Value = '''59.29'
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D11').Value = '''0.0891'
$ws.Range('E11').Value = '  +4.42%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '2.563.39'
$ws.Range('E13').Value = '  +2.92%  '
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').Value = '''22.05'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '''0.801'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').Value = '''5.60'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '2.246.44'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D19').Value = '41.770.51'
$ws.Range('E19').Value = '  +5.43%  '
$ws.Range('D20').Value = '''72.09'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('E21').Value = '  -2.16%  '
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').Value = '''250.83'
$ws.Range('E23').Value = '  +10.14%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''2.39'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''2.32'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').Value = '''9.71'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('D29').Value = '''166.99'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('D30').Value = '''19.96'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').Value = '''1.42'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '''5.02'
$ws.Range('E34').Value = '  +6.75%  '
$ws.Range('E35').Value = '  +3.65%  '
$ws.Range('D36').Value = '''0.0636'
$ws.Range('E36').Value = '  +3.29%  '
$ws.Range('E37').Value = '  -4.78%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('D39').Value = '''2.38'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  +30.53%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  +5.23%  '
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('D44').Value = '''8.59'
$ws.Range('E44').Value = '  +8.83%  '
$ws.Range('D45').Value = '''0.0981'
$ws.Range('E45').Value = '  +6.62%  '
$ws.Range('D46').Value = '''1.22'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '''99.00'
$ws.Range('D48').Value = '1.478.38'
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''16.49'
$ws.Range('E49').Value = '  -6.64%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.81'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '''52.56'
$ws.Range('E51').Value = '  +8.76%  '
